$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 185, shifting the existing rows 185-196 down to 186-197.
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with the new weekly price-point record.
$ws.Cells.Item(185, 1).Value = 7
$ws.Cells.Item(185, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(185, 3).Value = "Ñuble"
$ws.Cells.Item(185, 4).Value = 44585
$ws.Cells.Item(185, 5).Value = 16
$ws.Cells.Item(185, 6).Value = 100112043
$ws.Cells.Item(185, 7).Value = "Pepino ensalada"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 100
$ws.Cells.Item(185, 11).Value = 8000
$ws.Cells.Item(185, 12).Value = 8500
$ws.Cells.Item(185, 13).Value = 8250
$ws.Cells.Item(185, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(185, 15).Value = "Región del Maule"
$ws.Cells.Item(185, 16).Value = 103
$ws.Cells.Item(185, 17).Value = 80
$ws.Cells.Item(185, 18).Value = "Hortaliza"
